# Automatische test-sync: 2025-08-19 19:43:50
$wb = $excel.ActiveWorkbook

# --- Logs sheet: append a new row (row 8) mirroring the log entry added upstream ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Interne taak"
$logs.Range("B8").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("D8").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F8").Value = "2025-08-19 19:43:33"
$logs.Range("G8").Value = "Nee"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# --- Extend conditional formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`7")
    $newRange = $logs.Range("$col`2:$col`8")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the tally for "Intern verzoek / Actie voor medewerker" ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 7
